$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: replace values with 2-decimal rounded values (custom accuracy)
$ws.Range("B5").Value = 14.21
$ws.Range("C5").Value = 11.01
$ws.Range("D5").Value = 0.64
$ws.Range("E5").Value = 31.01
$ws.Range("F5").Value = 25.78
$ws.Range("G5").Value = 11.36
$ws.Range("H5").Value = 44.03
$ws.Range("I5").Value = 17.4
$ws.Range("J5").Value = 7.83
$ws.Range("K5").Value = 12.07
$ws.Range("L5").Value = 13.25
$ws.Range("M5").Value = 13.16
$ws.Range("N5").Value = 3.46
$ws.Range("O5").Value = 11.05
$ws.Range("P5").Value = 16.14
$ws.Range("Q5").Value = 9.32
$ws.Range("R5").Value = 0.11
$ws.Range("S5").Value = 0.39
$ws.Range("T5").Value = 164.18
$ws.Range("U5").Value = 31.51
$ws.Range("V5").Value = 10.12
$ws.Range("W5").Value = 21.2
$ws.Range("X5").Value = 11.49
$ws.Range("Y5").Value = 1.5
$ws.Range("Z5").Value = 21.81
$ws.Range("AA5").Value = 9.1
$ws.Range("AB5").Value = 8.58
$ws.Range("AC5").Value = 9.31
$ws.Range("AD5").Value = 13.57
$ws.Range("AE5").Value = 0.08
$ws.Range("AF5").Value = 39.95
$ws.Range("AG5").Value = 6.03
$ws.Range("AH5").Value = 12.88

# Remove row 6 entirely (data shrunk)
$ws.Rows(6).Delete()

# Column width tweaks: O (15) and Z (26) go from raw width 8 to raw width 7
# (ColumnWidth property is offset by ~0.83 from the raw OOXML width, so use 6.17)
$ws.Columns("O").ColumnWidth = 6.17
$ws.Columns("Z").ColumnWidth = 6.17
